$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix calculation of change: Initial biomass (B8) value reset to 0
$ws.Range("B8").Value = 0

# 2. Include final total biomass column in the change-in-biomass table
#    Header cell C12 gets a new label matching the style of B12
$ws.Range("B12").Copy()
$ws.Range("C12").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C12").Value = "Final total biomass" + [char]10 + "(tonnes CO2e)"
$ws.Range("C12").Font.Bold = $true
$ws.Range("C12").Characters(31, 1).Font.Subscript = $true
$ws.Range("C12").Characters(32, 2).Font.Bold = $true

# 3. Data row below the header now also has a value in column C (matching B13's style/value)
$ws.Range("B13").Copy()
$ws.Range("C13").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0

# 4. Update the active selection on the sheet
$ws.Range("B9").Select()
